$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values in the order that reproduces the shared-string table append order
# (25=Feb 15th, 26=Parentheses Not handled correctly, 27=Clion, cmake, clang, gcc,
#  28=Parentheses and Neg integers done, 29=Testing with different inputs,
#  30=Not correctly incrementing through string)
$ws.Range("A12").Value = "Feb 15th"
$ws.Range("E10").Value = "Parentheses Not handled correctly"
$ws.Range("F10").Value = "Clion, cmake, clang, gcc"
$ws.Range("C12").Value = "Parentheses and Neg integers done"
$ws.Range("D12").Value = "Testing with different inputs"
$ws.Range("E12").Value = "Not correctly incrementing through string"

# A10 reuses existing shared string "Feb 13th" (index 18), and F12 reuses
# "Clion, cmake, clang, gcc" (index 27), both already present by this point.
$ws.Range("A10").Value = "Feb 13th"
$ws.Range("F12").Value = "Clion, cmake, clang, gcc"

# Match the row heights used by the rest of the data rows
$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75

# Update selection to match new active cell
$ws.Range("D13").Select()
